# Append a new "Chmod 0777 ..." / forum-link row to Sheet1, then leave the
# selection where the author's session ended up (B13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 7: plain-text note in column A, URL text (styled like the other
# link cells but without an actual hyperlink) in column B.
$ws.Range("A7").Value = "Chmod 0777 on the folder that is being samba shared."
$ws.Range("B7").Value = "http://ubuntuforums.org/showthread.php?t=1723762"
$ws.Range("B7").Style = "Hyperlink"

# Final selection left on B13 (matches the saved sheetView selection).
$ws.Range("B13").Select()
